$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-11 from
# serial date 45233 (2023-11-03) to 45243 (2023-11-13).
$ws.Range("C2:C11").Value = 45243
